{"js": "// Update the date and the 25 division problems in the worksheet table.\nconst pairs = [\n  [\"2023-10-14 Saturday\", \"2023-10-15 Sunday\"],\n  [\"86\u00f73=\", \"53\u00f73=\"],\n  [\"74\u00f75=\", \"15\u00f78=\"],\n  [\"17\u00f75=\", \"62\u00f72=\"],\n  [\"26\u00f77=\", \"86\u00f78=\"],\n  [\"50\u00f76=\", \"10\u00f78=\"],\n  [\"11\u00f76=\", \"36\u00f73=\"],\n  [\"62\u00f78=\", \"67\u00f77=\"],\n  [\"59\u00f78=\", \"43\u00f78=\"],\n  [\"39\u00f78=\", \"67\u00f79=\"],\n  [\"44\u00f79=\", \"71\u00f79=\"],\n  [\"63\u00f77=\", \"70\u00f78=\"],\n  [\"56\u00f73=\", \"26\u00f72=\"],\n  [\"16\u00f79=\", \"11\u00f79=\"],\n  [\"12\u00f75=\", \"75\u00f75=\"],\n  [\"35\u00f78=\", \"45\u00f75=\"],\n  [\"10\u00f72=\", \"28\u00f72=\"],\n  [\"99\u00f76=\", \"74\u00f79=\"],\n  [\"66\u00f75=\", \"89\u00f76=\"],\n  [\"10\u00f76=\", \"67\u00f76=\"],\n  [\"39\u00f74=\", \"79\u00f75=\"],\n  [\"14\u00f74=\", \"88\u00f72=\"],\n  [\"89\u00f79=\", \"43\u00f79=\"],\n  [\"88\u00f79=\", \"63\u00f73=\"],\n  [\"94\u00f78=\", \"29\u00f79=\"],\n  [\"54\u00f79=\", \"52\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date and the 25 division problems in the worksheet table.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2023-10-14 Saturday\", \"2023-10-15 Sunday\"),\n    @(\"86\u00f73=\", \"53\u00f73=\"),\n    @(\"74\u00f75=\", \"15\u00f78=\"),\n    @(\"17\u00f75=\", \"62\u00f72=\"),\n    @(\"26\u00f77=\", \"86\u00f78=\"),\n    @(\"50\u00f76=\", \"10\u00f78=\"),\n    @(\"11\u00f76=\", \"36\u00f73=\"),\n    @(\"62\u00f78=\", \"67\u00f77=\"),\n    @(\"59\u00f78=\", \"43\u00f78=\"),\n    @(\"39\u00f78=\", \"67\u00f79=\"),\n    @(\"44\u00f79=\", \"71\u00f79=\"),\n    @(\"63\u00f77=\", \"70\u00f78=\"),\n    @(\"56\u00f73=\", \"26\u00f72=\"),\n    @(\"16\u00f79=\", \"11\u00f79=\"),\n    @(\"12\u00f75=\", \"75\u00f75=\"),\n    @(\"35\u00f78=\", \"45\u00f75=\"),\n    @(\"10\u00f72=\", \"28\u00f72=\"),\n    @(\"99\u00f76=\", \"74\u00f79=\"),\n    @(\"66\u00f75=\", \"89\u00f76=\"),\n    @(\"10\u00f76=\", \"67\u00f76=\"),\n    @(\"39\u00f74=\", \"79\u00f75=\"),\n    @(\"14\u00f74=\", \"88\u00f72=\"),\n    @(\"89\u00f79=\", \"43\u00f79=\"),\n    @(\"88\u00f79=\", \"63\u00f73=\"),\n    @(\"94\u00f78=\", \"29\u00f79=\"),\n    @(\"54\u00f79=\", \"52\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
